# Update countries & provincias Spain
# - Australia overtakes El Salvador and Costa Rica in the ranking (rows 74-76)
# - Islas Malvinas swaps places with Groenlandia (rows 210-211, tied values)
# - Refresh standalone country counters for Bolivia (row 36) and Vietnam (row 162)
# - Refresh the "Datos actualizados..." timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 05:09"

# --- Ranking shuffle: Australia moves above El Salvador and Costa Rica ---
$ws.Range("A74").Value = "Australia"
$ws.Range("B74").Value = 13948
$ws.Range("C74").Value = 353
$ws.Range("D74").Value = 8775
$ws.Range("E74").Value = 5028
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 6
$ws.Range("H74").Value = 145

$ws.Range("A75").Value = "El Salvador"
$ws.Range("B75").Value = 13792
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 7415
$ws.Range("E75").Value = 5998
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 379

$ws.Range("A76").Value = "Costa Rica"
$ws.Range("B76").Value = 13669
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 3505
$ws.Range("E76").Value = 10077
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 87

# --- Tied-value swap: Islas Malvinas now listed before Groenlandia ---
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

# --- Bolivia (row 36) updated counters ---
$ws.Range("B36").Value = 66456
$ws.Range("C36").Value = 1204
$ws.Range("D36").Value = 20614
$ws.Range("E36").Value = 43369
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 66
$ws.Range("H36").Value = 2473

# --- Vietnam (row 162) updated counters ---
$ws.Range("B162").Value = 415
$ws.Range("C162").Value = 2
$ws.Range("D162").Value = 365
$ws.Range("E162").Value = 50
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 0
